# Updates cryptos list data (Price and Volume(1h) columns) per latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.102.15"
$ws.Range("E2").Value = "  -2.75%  "

$ws.Range("D3").Value = "1.716.88"
$ws.Range("E3").Value = "  -2.99%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.31"
$ws.Range("E5").Value = "  -6.01%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4726"
$ws.Range("E7").Value = "  +5.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3435"
$ws.Range("E8").Value = "  -3.63%  "

$ws.Range("E9").Value = "  +0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07278"
$ws.Range("E10").Value = "  -2.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.044"

$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.93"
$ws.Range("E13").Value = "  -4.92%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.881"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").Value = "1.718.33"
$ws.Range("E15").Value = "  -3.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.892"
$ws.Range("E16").Value = "  -4.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.20"
$ws.Range("E17").Value = "  -4.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06359"
$ws.Range("E19").Value = "  -1.28%  "

$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.54"
$ws.Range("E21").Value = "  -3.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.631"
$ws.Range("E22").Value = "  -2.54%  "

$ws.Range("D23").Value = "27.135.84"
$ws.Range("E23").Value = "  -2.81%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.87"
$ws.Range("E24").Value = "  -3.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.123"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.68"
$ws.Range("E26").Value = "  -3.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.53"
$ws.Range("E27").Value = "  -4.08%  "

$ws.Range("D28").Value = "1.908.65"
$ws.Range("E28").Value = "  -3.34%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.097"
$ws.Range("E29").Value = "  -2.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.63"
$ws.Range("E30").Value = "  -4.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.018"
$ws.Range("E31").Value = "  -7.87%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09168"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.593"
$ws.Range("E33").Value = "  -1.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.331"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02207"
$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05834"
$ws.Range("E36").Value = "  -4.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.01"
$ws.Range("E37").Value = "  -7.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2000"
$ws.Range("E38").Value = "  -4.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.745"
$ws.Range("E39").Value = "  -4.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.393"
$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5901"
$ws.Range("E41").Value = "  -6.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.122"
$ws.Range("E42").Value = "  -5.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.471"
$ws.Range("E43").Value = "  -5.29%  "

$ws.Range("E44").Value = "  -5.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5665"
$ws.Range("E45").Value = "  -4.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.562"
$ws.Range("E46").Value = "  -4.76%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.67"
$ws.Range("E47").Value = "  -3.77%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.844"
$ws.Range("E48").Value = "  -5.76%  "

$ws.Range("E49").Value = "  -3.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.087"
$ws.Range("E50").Value = "  -4.43%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  +0.09%  "
